# ShantanuMane Cover Letter (Respawn) - game-related content rewrite
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the title paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Drop the trailing space after "...involved in them. " and, right
#    after it, insert the two new paragraphs about combat/animation
#    passion and the Titanfall power fantasy (these land right before
#    the "I love Action Games..." paragraph).
# ---------------------------------------------------------------------
$newParas = "collaborating with and learning from the people involved in them.^p" + `
"Combat, weapons, action and animation are where my passion truly lies. I have taken and continue to take time to learn combat design and combat systems. I am also putting time into learning about animation programming and am working towards creating a gameplay animation system related to but not limited to combat. ^p" + `
"The free-flowing movement of Titanfall fulfills the power fantasy of being a nimble and agile swashbuckler as the Pilot and then the powerful arsenal of weapons you have when using your Titan make you go gung-ho, keep the trigger pulled and lay all your firepower into your enemies. I think the game does an excellent job of making difficult feats achievable with deftness. This puts the player in the power fantasy the game wants them to feel from the very first moment they wall-run and leaves room for creativity past that point, with the systems being intuitive and rewarding the player for engaging in them."

$d.Content.Find.Execute(
    "collaborating with and learning from the people involved in them. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newParas,
    2) | Out-Null

# ---------------------------------------------------------------------
# 3. Rewrite the back half of the "I love Action Games..." paragraph -
#    swap the Uncharted references for Titanfall 2 impressions.
# ---------------------------------------------------------------------
$oldUncharted = " The Uncharted games nail down these feelings masterfully and immaculately. The action, the firefights and the chase sequences in them have made me ride waves of excitement and adrenaline. The high stakes and high adrenaline climactic chase sequences are particularly my favorites of these. They have had me on the edge of my seat, tightly gripping my controller, being on point with my shooting and pushing hard on the movement stick to charge through to make it out alive if only by the skin of my teeth! I think they are just very beautifully done."

$newTitanfall2 = " I’ve played Titanfall 2 and I can say it evoked these feelings in me gracefully and masterfully. I like parkour and I loved wall-running in the game. It being so fluid and natural to pull-off opened up avenues for me to explore all my abilities and get creative with them. I’ve had adrenaline-filled combat encounters in the game where I was on point, pushing my abilities to their extent, wall-running and shooting enemies then launching off and landing to go sliding into an enemy to finish them off with a powerful and bone-shattering melee strike!"

$d.Content.Find.Execute($oldUncharted, $true, $false, $false, $false, $false, $true, 1, $false, $newTitanfall2, 2) | Out-Null

# ---------------------------------------------------------------------
# 4. Shrink the old "Combat, weapons..." paragraph down to the short
#    editorial note that remains, then drop a fresh "_GoBack" bookmark
#    into the middle of its text.
# ---------------------------------------------------------------------
$oldCombatIntro = "Combat, weapons, action and animation are where my passion truly lies. I have and continue to take time to learn combat design and combat systems. I have also taken various design classes. I am also putting time into learning about animation programming and am working towards creating a gameplay animation system related to but not limited to combat. My knowledge of these and experience playing Uncharted 3 & 4 clearly told me how much the melee system had improved. And looking at the combat and animation in The Last of Us Part II I think the studio is making great headway. The action shown in the E3 2018 gameplay footage was nothing short of sublime!"

$newTalkAbout = "Talk about wanting to create power fantasy that they go for here? -"

$d.Content.Find.Execute($oldCombatIntro, $true, $false, $false, $false, $false, $true, 1, $false, $newTalkAbout, 2) | Out-Null

# Re-home the "_GoBack" bookmark between "go for " and "here? -"
$found = $d.Content.Find.Execute("go for here? -", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $anchor = $d.Content.Find.Parent
}
$markerText = "go for "
$searchRange = $d.Content
$searchRange.Find.Execute($markerText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertAt = $searchRange.End
$bmRange = $d.Range($insertAt, $insertAt)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 5. Rewrite the opening of the "Games with Action and..." paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Games with Action and engaging, compelling experiences are what ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Games with fluid, fast-paced action are what ",
    2) | Out-Null
